# Commit: "falei que sou inteligente"
#
# 1. The second paragraph ("Domingo o Grêmio ganhou o GREnal") gets "GREnal"
#    broken out into its own run, bracketed by proofErr spell-check markers
#    (as Word's background spell-checker does for words it doesn't recognise).
# 2. A new blank paragraph is added.
# 3. A new paragraph "A diele é super inteligente " is added (with proofErr
#    markers around "diele" and "super"), and the trailing _GoBack bookmark
#    now sits at the end of this new last paragraph.

$d = $word.ActiveDocument

$para2 = $d.Paragraphs.Item(2)
$p2Start = $para2.Range.Start
$p2End = $para2.Range.End - 1   # exclude the paragraph mark itself

# Wipe the old run text but keep the paragraph (and its trailing bookmark)
# intact so the bookmark re-anchors naturally after the new content lands.
$oldTextRange = $d.Range($p2Start, $p2End)
$oldTextRange.Text = ""

# Insert the replacement content (3 paragraphs worth of markup) at the spot
# the old text used to start, i.e. right before the bookmark.
$insertionPoint = $d.Range($p2Start, $p2Start)

$newXml = @"
<?xml version="1.0"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:r><w:t xml:space="preserve">Domingo o Gr&#234;mio ganhou o </w:t></w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r><w:t>GREnal</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
          </w:p>
          <w:p/>
          <w:p>
            <w:r><w:t xml:space="preserve">A </w:t></w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r><w:t>diele</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r><w:t xml:space="preserve"> &#233; </w:t></w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r><w:t>super</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r><w:t xml:space="preserve"> inteligente </w:t></w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@

$insertionPoint.InsertXML($newXml)

$d.Save()
